$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Reposition the background rectangle ("Retângulo 38", id=39) on the
#    "Máquinas" dashboard slide (sldId 260 -> the 4th slide in the deck).
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(4)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Id -eq 39) {
        $sh.Left = -6.158976377952756
        $sh.Top = 9.389922
    }
}

# ---------------------------------------------------------------------------
# 2) Bump the cached "datetimeFigureOut" footer date from 24/10/2020 to
#    25/10/2020 on every slide layout (Header & Footer date placeholder).
# ---------------------------------------------------------------------------
$oldDate = "24/10/2020"
$newDate = "25/10/2020"

$master = $p.SlideMaster
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    for ($si = 1; $si -le $layout.Shapes.Count; $si++) {
        $shp = $layout.Shapes.Item($si)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}
